# Generate Report for Handoff
# Updates the localization-status report: flips the "In Translation" status
# to "Ready for handoff" on every sheet, refreshes the associated timestamps,
# and widens the Status column(s) to fit the new, longer text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet (row 2: zh-cn status, de-de status, latest HO xliff date) ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-27 20:57:28"

# --- zh-cn detail sheet (Status + Latest Handoff Datetime) ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-27 20:57:23"

# --- de-de detail sheet (Status + Latest Handoff Datetime) ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-27 20:57:28"

# --- Widen the Status columns so the longer "Ready for handoff" text fits ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332
$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333332
$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333332
